$wb = $excel.ActiveWorkbook

# Update the shared string "Ready for handoff" -> "In Translation"
# This text appears in:
#  - Overview sheet: E2 (zh-cn status), F2 (de-de status)
#  - zh-cn sheet: C2 (Status)
#  - de-de sheet: C2 (Status)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Adjust column widths to reflect the new (shorter) text extent
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
